# Price update for 2026-02-07: append a new tracked-price row to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 38
$scratch = "Z100"

function Set-TextValue {
    param($cellAddress, [string]$text)

    # Writing a quoted string literal through a formula keeps the result as
    # genuine text (Excel's smart "Value" parser never gets a chance to
    # reinterpret it as a date/number). Copying that computed value and
    # pasting it back as a literal (PasteSpecial values-only, paste type
    # -4163 / xlPasteValues) converts it into a plain shared-string cell
    # with default styling - exactly like the source cells already on the
    # sheet - instead of leaving a formula or forcing a new number format.
    $escaped = $text.Replace('"', '""')
    $ws.Range($scratch).Formula = '="' + $escaped + '"'
    $ws.Range($scratch).Copy()
    $ws.Range($cellAddress).PasteSpecial(-4163)
}

Set-TextValue "A$targetRow" "2026-02-07"
Set-TextValue "B$targetRow" "175600"
Set-TextValue "C$targetRow" "40"
Set-TextValue "D$targetRow" "0"

$ws.Range($scratch).Clear()
$excel.CutCopyMode = 0

Write-Host "A$targetRow -> $($ws.Range("A$targetRow").Text)"
Write-Host "B$targetRow -> $($ws.Range("B$targetRow").Text)"
Write-Host "C$targetRow -> $($ws.Range("C$targetRow").Text)"
Write-Host "D$targetRow -> $($ws.Range("D$targetRow").Text)"
